$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    12 = @(0.1644653861563335, 0.004166984081917661, -0.009744353085738667, 1.829707928933755)
    13 = @(0.1568772835360597, 0.003563904471037351, 0.3221801397222507, 2.100310050176578)
    14 = @(0.1717729212645522, 0.004079320509223575, -0.03643626733293952, 1.679618112128374)
    15 = @(0.2268026515618569, 0.02166618773604188, 1.384061321255182, 3.823793060073173)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
    $ws.Range("K$row").Value = $values[2]
    $ws.Range("L$row").Value = $values[3]
}
